$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.562.10"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.812.31"
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'228.70"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  +8.59%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'36.64"
$ws.Range("E8").Value = "  +4.95%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.0701"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "'0.0966"
$ws.Range("E11").Value = "  +1.61%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "1.823.83"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("E15").Value = "  +1.64%  "
$ws.Range("E16").Value = "  +3.90%  "
$ws.Range("D17").Value = "34.540.63"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "'70.27"
$ws.Range("E18").Value = "  +1.50%  "
$ws.Range("D19").Value = "'246.94"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "'11.65"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("D23").Value = "'4.23"
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("D24").Value = "'2.25"
$ws.Range("E24").Value = "  +7.57%  "
$ws.Range("D25").Value = "'173.02"
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("D26").Value = "'8.03"
$ws.Range("E26").Value = "  +7.55%  "
$ws.Range("D27").Value = "'17.32"
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("E30").Value = "  +0.15%  "
$ws.Range("D31").Value = "'3.87"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "1.404.45"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").Value = "'2.46"
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("D40").Value = "'0.970"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").Value = "'82.93"
$ws.Range("E41").Value = "  -4.00%  "
$ws.Range("D42").Value = "'2.85"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "'1.20"
$ws.Range("E44").Value = "  +7.83%  "
$ws.Range("D45").Value = "'13.53"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("B47").Value = "Kaspa"
$ws.Range("C47").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D47").Value = "'0.0496"
$ws.Range("E47").Value = "  -5.30%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.972.46"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").Value = "'104.51"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -2.90%  "
